# fix: adjust test order and teardown
#
# The LoginData sheet drives the login test with a freshly generated,
# not-yet-registered email address in A2. Re-ordering the tests changed
# which generated address lands there by the time this sheet is read, and
# the teardown step (which used to leave several stale addresses behind
# in earlier runs) now properly cleans those up so only the one actually
# in use remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")
$ws.Activate()

# Email address generated earlier in the (reordered) run - no longer the
# one the login test ends up using, but it passes through A2 first.
$ws.Range("A2").Value = "juan.perez99_20017@test.com"

# Final email address left in place for the login test after teardown.
$ws.Range("A2").Value = "juan.perez99_25911@test.com"

# Reflect the cursor position left behind by the adjusted test/teardown
# order.
$ws.Range("B2").Select()
